$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Sending=FAPs, Ligand=Il13, Receptor=Il13ra2, Target=ECs (labels unchanged, TPM values refreshed) ---
$ws.Range("I2").Value = 0.6540874079906115
$ws.Range("J2").Value = 0.7393359457808691
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.0568385
$ws.Range("N2").Value = 0.113677
$ws.Range("O2").Value = 0.01274651757362603
$ws.Range("P2").Value = 0.008533937711420974
$ws.Range("Q2").Value = 0.006673730369833332
$ws.Range("R2").Value = 0.040042382219
$ws.Range("S2").Value = 0.00833733664063983
$ws.Range("T2").Value = 0.006309446909108452

# --- Row 3: Sending=FAPs, Ligand=Il13, Receptor=Il13ra2, Target=FAPs (labels unchanged, TPM values refreshed) ---
$ws.Range("I3").Value = 0.6540874079906115
$ws.Range("J3").Value = 0.7393359457808691
$ws.Range("O3").Value = 0.9872534824263741
$ws.Range("P3").Value = 0.9914660622885791
$ws.Range("S3").Value = 0.6457500713499718
$ws.Range("T3").Value = 0.7330264988717607

# --- Row 4 (new): Sending=MuSCs, Ligand=Il13, Receptor=Il13ra2, Target=ECs ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Il13"
$ws.Range("C4").Value = "Il13ra2"
$ws.Range("D4").Value = "ECs"

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.062095
$ws.Range("H4").Value = 0.12419
$ws.Range("I4").Value = 0.3459125920093885
$ws.Range("J4").Value = 0.2606640542191307
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.0568385
$ws.Range("N4").Value = 0.113677
$ws.Range("O4").Value = 0.01274651757362603
$ws.Range("P4").Value = 0.008533937711420974
$ws.Range("Q4").Value = 0.0035293866575
$ws.Range("R4").Value = 0.01411754663
$ws.Range("S4").Value = 0.004409180932986203
$ws.Range("T4").Value = 0.002224490802312522

# --- Row 5 (new): Sending=MuSCs, Ligand=Il13, Receptor=Il13ra2, Target=FAPs ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Il13"
$ws.Range("C5").Value = "Il13ra2"
$ws.Range("D5").Value = "FAPs"

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.062095
$ws.Range("H5").Value = 0.12419
$ws.Range("I5").Value = 0.3459125920093885
$ws.Range("J5").Value = 0.2606640542191307
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.402301
$ws.Range("N5").Value = 13.206903
$ws.Range("O5").Value = 0.9872534824263741
$ws.Range("P5").Value = 0.9914660622885791
$ws.Range("Q5").Value = 0.273360880595
$ws.Range("R5").Value = 1.64016528357
$ws.Range("S5").Value = 0.3415034110764024
$ws.Range("T5").Value = 0.2584395634168182

Write-Host "edit applied"
